$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column E (Obrigatorio) from "N" to "S" for rows 2 through 8
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 5).Value = "S"
}

$wb.Save()
